# Update giao diện tổng quan và chức năng
#
# 1. Sheet "Tiêu chuẩn TMĐT ĐG (CTN009)" (1st sheet):
#    - formulas in B8:D11 change from ROUNDUP(x/1.08,0) to ROUNDUP(x*1.08,0)
#    - becomes the active / selected sheet, selection moves to F11
# 2. Sheet "Tiêu chuẩn TMĐT(CTN007)" (2nd sheet):
#    - formulas in B9:F13 change from ROUNDUP(x/1.08,0) to ROUNDUP(x*1.08,0)
#    - is no longer the active sheet, selection moves to H11

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 2 (Tiêu chuẩn TMĐT(CTN007)): flip formulas from /1.08 to *1.08 ---
$ws2.Range("B9").Formula  = "=ROUNDUP(B2*1.08,0)"
$ws2.Range("C9").Formula  = "=ROUNDUP(C2*1.08,0)"
$ws2.Range("D9").Formula  = "=ROUNDUP(D2*1.08,0)"
$ws2.Range("E9").Formula  = "=ROUNDUP(E2*1.08,0)"
$ws2.Range("F9").Formula  = "=ROUNDUP(F2*1.08,0)"

$ws2.Range("B10").Formula = "=ROUNDUP(B3*1.08,0)"
$ws2.Range("C10").Formula = "=ROUNDUP(C3*1.08,0)"
$ws2.Range("D10").Formula = "=ROUNDUP(D3*1.08,0)"
$ws2.Range("E10").Formula = "=ROUNDUP(E3*1.08,0)"
$ws2.Range("F10").Formula = "=ROUNDUP(F3*1.08,0)"

$ws2.Range("B11").Formula = "=ROUNDUP(B4*1.08,0)"
$ws2.Range("C11").Formula = "=ROUNDUP(C4*1.08,0)"
$ws2.Range("D11").Formula = "=ROUNDUP(D4*1.08,0)"
$ws2.Range("E11").Formula = "=ROUNDUP(E4*1.08,0)"
$ws2.Range("F11").Formula = "=ROUNDUP(F4*1.08,0)"

$ws2.Range("B12").Formula = "=ROUNDUP(B5*1.08,0)"
$ws2.Range("C12").Formula = "=ROUNDUP(C5*1.08,0)"
$ws2.Range("D12").Formula = "=ROUNDUP(D5*1.08,0)"
$ws2.Range("E12").Formula = "=ROUNDUP(E5*1.08,0)"
$ws2.Range("F12").Formula = "=ROUNDUP(F5*1.08,0)"

$ws2.Range("B13").Formula = "=ROUNDUP(B6*1.08,0)"
$ws2.Range("C13").Formula = "=ROUNDUP(C6*1.08,0)"
$ws2.Range("D13").Formula = "=ROUNDUP(D6*1.08,0)"
$ws2.Range("E13").Formula = "=ROUNDUP(E6*1.08,0)"
$ws2.Range("F13").Formula = "=ROUNDUP(F6*1.08,0)"

# --- Sheet 1 (Tiêu chuẩn TMĐT ĐG (CTN009)): flip formulas from /1.08 to *1.08 ---
$ws1.Range("B8").Formula  = "=ROUNDUP(B2*1.08,0)"
$ws1.Range("C8").Formula  = "=ROUNDUP(C2*1.08,0)"
$ws1.Range("D8").Formula  = "=ROUNDUP(D2*1.08,0)"

$ws1.Range("B9").Formula  = "=ROUNDUP(B3*1.08,0)"
$ws1.Range("C9").Formula  = "=ROUNDUP(C3*1.08,0)"
$ws1.Range("D9").Formula  = "=ROUNDUP(D3*1.08,0)"

$ws1.Range("B10").Formula = "=ROUNDUP(B4*1.08,0)"
$ws1.Range("C10").Formula = "=ROUNDUP(C4*1.08,0)"
$ws1.Range("D10").Formula = "=ROUNDUP(D4*1.08,0)"

$ws1.Range("B11").Formula = "=ROUNDUP(B5*1.08,0)"
$ws1.Range("C11").Formula = "=ROUNDUP(C5*1.08,0)"
$ws1.Range("D11").Formula = "=ROUNDUP(D5*1.08,0)"

# --- Selections / active sheet ---
# Select sheet2's new cell first (it becomes the non-active selection),
# then select sheet1's new cell last so sheet1 ends up the active /
# tab-selected sheet, matching the diff (activeTab moves from 1 -> 0).
$ws2.Range("H11").Select()
$ws1.Range("F11").Select()
